$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet is protected; unprotect so the cells can be edited, then re-protect afterwards.
$ws.Unprotect()

# Update the confidential disclosure date string (2021-05-27 -> 2021-05-28)
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."
# Re-apply auto row height so the multi-line text doesn't leave a stray explicit row height.
$ws.Rows(9).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2532974909413183
$ws.Range("E2").Value = -0.001030337721808716

$ws.Range("D3").Value = 0.2537762012780171
$ws.Range("E3").Value = 0.001845991561181481

$ws.Range("D4").Value = 0.2443829449754355
$ws.Range("E4").Value = 0

$ws.Range("D5").Value = 0.2485433628052291
$ws.Range("E5").Value = -0.003289057558507169

$ws.Range("D6").Value = 0.9999999999999999
$ws.Range("E6").Value = -0.0006099866598198167

# Restore sheet protection to its original (protected) state.
$ws.Protect()
